$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F (reuse the same header formatting as the other headers)
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Replace the text timestamps in column A with real date/time serial values
$ws.Range("A2").Value = 45685.64807152778
$ws.Range("A3").Value = 45685.64943032408
$ws.Range("A4").Value = 45685.65091064815
$ws.Range("A5").Value = 45685.64806805555
$ws.Range("A6").Value = 45685.64942685185
$ws.Range("A7").Value = 45685.65090717593

# Fix up row 5 values (previously a different reading) and add Trening column
$ws.Range("B5").Value = 1108
$ws.Range("C5").Value = 8.65
$ws.Range("D5").Value = 3.391879796981811

# Trening labels for the existing rows (part 1 of training session)
$ws.Range("F2").Value = "Duża Gra"
$ws.Range("F3").Value = "Duża Gra"
$ws.Range("F4").Value = "Duża Gra"
$ws.Range("F5").Value = "Duża Gra"
$ws.Range("F6").Value = "Duża Gra"
$ws.Range("F7").Value = "Duża Gra"

# New rows 8-13 (part 2 of training session)
$ws.Range("A8").Value = 45685.6674199074
$ws.Range("B8").Value = 2780
$ws.Range("C8").Value = 14.06
$ws.Range("D8").Value = 3.664245332990375
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Mała Gra"

$ws.Range("A9").Value = 45685.67531342593
$ws.Range("B9").Value = 3462
$ws.Range("C9").Value = 14.61
$ws.Range("D9").Value = 3.076576471328735
$ws.Range("E9").Value = "10-15"
$ws.Range("F9").Value = "Mała Gra"

$ws.Range("A10").Value = 45685.68058773148
$ws.Range("B10").Value = 3917.7
$ws.Range("C10").Value = 14.27
$ws.Range("D10").Value = 3.130322422300065
$ws.Range("E10").Value = "10-15"
$ws.Range("F10").Value = "Mała Gra"

$ws.Range("A11").Value = 45685.66741643519
$ws.Range("B11").Value = 2779.7
$ws.Range("C11").Value = 8.92
$ws.Range("D11").Value = 2.893186858722143
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "Mała Gra"

$ws.Range("A12").Value = 45685.66908888889
$ws.Range("B12").Value = 2924.2
$ws.Range("C12").Value = 9.32
$ws.Range("D12").Value = 2.822287797927857
$ws.Range("E12").Value = "5-10"
$ws.Range("F12").Value = "Mała Gra"

$ws.Range("A13").Value = 45685.68470462963
$ws.Range("B13").Value = 4273.4
$ws.Range("C13").Value = 9.57
$ws.Range("D13").Value = 2.978939294815061
$ws.Range("E13").Value = "5-10"
$ws.Range("F13").Value = "Mała Gra"

# Apply the date/time display format to the whole Timestamp column (A2:A13).
# A lowercase format is tried on the first cell, then the final uppercase
# format is applied to the whole column - this mirrors the source workbook,
# which ends up with two numFmt entries (164 unused, 165 applied) while
# keeping a single style slot for all the timestamp cells.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
